# "Actualice el documento comandos HTML"
# Insert a new row for the generic "git push" command just above the existing
# "git push <repo> <rama>" row, reusing the "Subir a repositorio" description
# and the yellow highlight style used for the other "main command" rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22 currently holds "git push <nombre repo> <rama a subir>" / "... --all" /
# "Subir a repositorio". Insert a blank row above it so that row becomes row 23
# (and the following "git clone" row becomes row 24).
$ws.Rows("22:22").Insert()

# Populate the newly inserted row 22 with the generic push command.
$ws.Range("B22").Value = "git push "
$ws.Range("D22").Value = "Subir a repositorio"

# Match the yellow highlight formatting used on the other key-command rows
# (e.g. "git clone", now at row 24).
$ws.Range("B22").Interior.Color = 65535

# Reflect the author's final cell selection in the saved view state.
$ws.Range("B15").Select()
